$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "AAA" value in C2 gained trailing whitespace in the source test data,
# used to exercise the new AutomaticallyTrimAllStringValues reader option.
$ws.Range("C2").Value = "AAA   "
